# Insert one new week of data (3 rows) for "Comercializadora del Agro de Limarí - Limón"
# at the top of the historical log (row 991), pushing all later rows down by 3.
# New data corresponds to the week of 2023-11-09 (serial 45239), quality grades
# 1a/2a/3a amarillo, packed in "$/malla 18 kilos".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows starting at row 991 (existing data shifts down).
$ws.Range("A991:T993").EntireRow.Insert()

# Copy the date number format used by the rest of column D down into the new rows
# (row 994 is the row that used to be 991, already carrying the correct format).
$ws.Range("D991:D993").NumberFormat = $ws.Range("D994").NumberFormat

$rows = @(
    @{ Row = 991; Date = 45239; Calidad = "1a amarillo"; Volumen = 600; PMin = 7800; PMax = 8000; PProm = 7900; Unidad = "`$/malla 18 kilos"; PKg = 439; KgUnidad = 18 },
    @{ Row = 992; Date = 45239; Calidad = "2a amarillo"; Volumen = 560; PMin = 6800; PMax = 7000; PProm = 6900; Unidad = "`$/malla 18 kilos"; PKg = 383; KgUnidad = 18 },
    @{ Row = 993; Date = 45239; Calidad = "3a amarillo"; Volumen = 420; PMin = 4800; PMax = 5000; PProm = 4900; Unidad = "`$/malla 18 kilos"; PKg = 272; KgUnidad = 18 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = $r.Date
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100102
    $ws.Cells.Item($row, 8).Value = "Cítricos"
    $ws.Cells.Item($row, 9).Value = 100102003
    $ws.Cells.Item($row, 10).Value = "Limón"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
